$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 192 ("「きょうりゅうはどうやってがっこうにいくの？」" entry),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(192).Delete()
